# Add two "rerun" sheets to the NAS-PB benchmarking workbook:
#   - GNUPlot-Time-Class-A-Rerun  (copied from GNUPlot-Time-Class-B, then re-valued)
#   - GNUPlot-Time-Class-B-Rerun  (copied from GNUPlot-Time-Class-B, then re-valued)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "GNUPlot-Time-Class-A" (sheet10): the active cell moved back to A1
# ---------------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("GNUPlot-Time-Class-A")
$wsA.Activate()
$wsA.Range("A1").Select()

# ---------------------------------------------------------------------------
# Create "GNUPlot-Time-Class-A-Rerun" by duplicating "GNUPlot-Time-Class-B"
# (it carries over the same column layout / headers / styles as Class B)
# ---------------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("GNUPlot-Time-Class-B")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsB.Copy($null, $lastSheet)
$rerunA = $wb.Worksheets.Item($wb.Worksheets.Count)
$rerunA.Name = "GNUPlot-Time-Class-A-Rerun"

# Row 2
$rerunA.Cells.Item(2,2).Value = 1.18
$rerunA.Cells.Item(2,3).Value = 1.05
$rerunA.Cells.Item(2,4).Value = 1.34
$rerunA.Cells.Item(2,5).Value = 1.46
$rerunA.Cells.Item(2,6).Value = 78.66

# Row 3 -> cleared out (keep formatting, drop values)
$rerunA.Range("B3:F3").ClearContents()

# Row 4 - B4 used to hold the text placeholder "x"; it now gets a real number,
# so first copy the numeric formatting from its neighbour C4, then set value.
$rerunA.Cells.Item(4,3).Copy()
$rerunA.Cells.Item(4,2).PasteSpecial(-4122)
$excel.CutCopyMode = 0
$rerunA.Cells.Item(4,2).Value = 4.81
$rerunA.Cells.Item(4,3).Value = 4.97
$rerunA.Cells.Item(4,4).Value = 6.03
$rerunA.Cells.Item(4,5).Value = 4.49
$rerunA.Cells.Item(4,6).Value = 9.21

# Row 5
$rerunA.Cells.Item(5,2).Value = 0.64
$rerunA.Cells.Item(5,3).Value = 1.14
$rerunA.Cells.Item(5,4).Value = 0.91
$rerunA.Cells.Item(5,5).Value = 0.89
$rerunA.Cells.Item(5,6).Value = 6.76

# Row 6
$rerunA.Cells.Item(6,2).Value = 57.13
$rerunA.Cells.Item(6,3).Value = 32.13
$rerunA.Cells.Item(6,4).Value = 20.19
$rerunA.Cells.Item(6,5).Value = 16.79
$rerunA.Cells.Item(6,6).Value = 95.05

# Row 7 -> cleared out
$rerunA.Range("B7:F7").ClearContents()

$rerunA.Range("G9").Select()

# ---------------------------------------------------------------------------
# Create "GNUPlot-Time-Class-B-Rerun" by duplicating "GNUPlot-Time-Class-B"
# ---------------------------------------------------------------------------
$wsB2 = $wb.Worksheets.Item("GNUPlot-Time-Class-B")
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsB2.Copy($null, $lastSheet2)
$rerunB = $wb.Worksheets.Item($wb.Worksheets.Count)
$rerunB.Name = "GNUPlot-Time-Class-B-Rerun"

# Row 2
$rerunB.Cells.Item(2,2).Value = 80.63
$rerunB.Cells.Item(2,3).Value = 42.1
$rerunB.Cells.Item(2,4).Value = 42.94
$rerunB.Cells.Item(2,5).Value = 46.17
$rerunB.Cells.Item(2,6).Value = 159.74

# Row 3 -> cleared out
$rerunB.Range("B3:F3").ClearContents()

# Row 4 - B4 keeps the "x" placeholder text; C4:F4 get new numbers but drop
# down to the plain/default style (s=0) instead of the numeric style (s=3).
$rerunB.Range("C4:F4").ClearContents()
$rerunB.Range("C4:F4").Style = "Normal"
$rerunB.Cells.Item(4,3).Value = 63.81
$rerunB.Cells.Item(4,4).Value = 69.57
$rerunB.Cells.Item(4,5).Value = 52.4
$rerunB.Cells.Item(4,6).Value = 328.88

# Row 5 - B5,C5,D5,F5 stay on the numeric style; E5 drops to the default style
$rerunB.Cells.Item(5,2).Value = 2.86
$rerunB.Cells.Item(5,3).Value = 4.48
$rerunB.Cells.Item(5,4).Value = 3.57
$rerunB.Range("E5").ClearContents()
$rerunB.Range("E5").Style = "Normal"
$rerunB.Cells.Item(5,5).Value = 3.52
$rerunB.Cells.Item(5,6).Value = 32.94

# Row 6
$rerunB.Cells.Item(6,2).Value = 267.07
$rerunB.Cells.Item(6,3).Value = 140.01
$rerunB.Cells.Item(6,4).Value = 102.48
$rerunB.Cells.Item(6,5).Value = 74.21
$rerunB.Cells.Item(6,6).Value = 532.4

# Row 7 -> cleared out
$rerunB.Range("B7:F7").ClearContents()

$rerunB.Range("F12").Select()
$rerunB.Activate()
